$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 73, pushing existing rows 73:99 down to 74:100
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new record
$ws.Range("A73").Value = 5
$ws.Range("B73").Value = "Macroferia Regional de Talca"
$ws.Range("C73").Value = "Maule"
$ws.Range("D73").Value = 44855
$ws.Range("D73").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E73").Value = 7
$ws.Range("F73").Value = 100112026
$ws.Range("G73").Value = "Haba"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 400
$ws.Range("K73").Value = 7000
$ws.Range("L73").Value = 7000
$ws.Range("M73").Value = 7000
$ws.Range("N73").Value = "`$/saco 25 kilos"
$ws.Range("O73").Value = "Región de O'Higgins"
$ws.Range("P73").Value = 280
$ws.Range("Q73").Value = 25
$ws.Range("R73").Value = "Hortaliza"
